# Add an "ExportTemplate" header block above the transactions table and
# fix the "Category ID" column header typo to "Categorie ID".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows at the top of the sheet. This pushes the existing
# header row (and all data rows below it) down by three rows, carrying
# their styles/values along automatically.
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# Row 1: user name banner. Row 2: export timestamp banner. Row 3 stays blank
# as a spacer above the (now shifted) table header on row 4.
$ws.Range("A1").Value = "User Name: Kitsapas Chanlee"
$ws.Range("A2").Value = "Created At: 2024-07-03 10:43:08"

# Fix the column header typo: "Category ID" -> "Categorie ID".
# After the row insert, the header row lives at row 4; "Category ID" is
# the second header column (column B).
$ws.Range("B4").Value = "Categorie ID"
